$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 (the extra/duplicate "write_data" row), shifting rows 16-23 up
$ws.Rows.Item(15).Delete()

# Update selection to match final state
$ws.Range("B16").Select()
